$wb = $excel.ActiveWorkbook

# ---- Step 1: insert the new '2022-Q1' detail sheet right before the '总计' summary sheet ----
$srcTemplate = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Copy the 28-row block (header + 27 data rows) from the 2021-Q4 sheet so the new
# sheet inherits the same header/index-column formatting (bold+border header row,
# bold+border column A) without fabricating new style entries.
$srcTemplate.Range("A1:H28").Copy($q1.Range("A1"))

# ---- Step 2: overwrite header text for the new sheet ----
$q1.Range("B1").Value2 = '基金代码'
$q1.Range("C1").Value2 = '基金名称'
$q1.Range("D1").Value2 = '基金规模'
$q1.Range("E1").Value2 = '股票总仓位'
$q1.Range("F1").Value2 = '仓位占比'
$q1.Range("G1").Value2 = '持有市值(亿元)'
$q1.Range("H1").Value2 = '仓位排名'

# ---- Step 3: fill in the 27 fund-holding data rows ----
# row 2
$q1.Range("A2").Value2 = 0
$q1.Range("B2").Value2 = '''001955'
$q1.Range("C2").Value2 = '中欧养老产业混合'
$q1.Range("D2").Value2 = '''46.95'
$q1.Range("E2").Value2 = '''92.35'
$q1.Range("F2").Value2 = '''8.71'
$q1.Range("G2").Value2 = '''4.0893'
$q1.Range("H2").Value2 = 7

# row 3
$q1.Range("A3").Value2 = 1
$q1.Range("B3").Value2 = '''001071'
$q1.Range("C3").Value2 = '华安媒体互联网混合'
$q1.Range("D3").Value2 = '''51.61'
$q1.Range("E3").Value2 = '''92.88'
$q1.Range("F3").Value2 = '''4.32'
$q1.Range("G3").Value2 = '''2.2296'
$q1.Range("H3").Value2 = 7

# row 4
$q1.Range("A4").Value2 = 2
$q1.Range("B4").Value2 = '''166027'
$q1.Range("C4").Value2 = '中欧创业板两年定期开放混合A'
$q1.Range("D4").Value2 = '''21.11'
$q1.Range("E4").Value2 = '''99.64'
$q1.Range("F4").Value2 = '''9.93'
$q1.Range("G4").Value2 = '''2.0962'
$q1.Range("H4").Value2 = 1

# row 5
$q1.Range("A5").Value2 = 3
$q1.Range("B5").Value2 = '''010429'
$q1.Range("C5").Value2 = '中欧睿见混合'
$q1.Range("D5").Value2 = '''27.61'
$q1.Range("E5").Value2 = '''92.03'
$q1.Range("F5").Value2 = '''7.41'
$q1.Range("G5").Value2 = '''2.0459'
$q1.Range("H5").Value2 = 7

# row 6
$q1.Range("A6").Value2 = 4
$q1.Range("B6").Value2 = '''159766'
$q1.Range("C6").Value2 = '富国中证旅游主题交易型开放式指数证券投资基金'
$q1.Range("D6").Value2 = '''17.47'
$q1.Range("E6").Value2 = '''99.35'
$q1.Range("F6").Value2 = '''10.72'
$q1.Range("G6").Value2 = '''1.8728'
$q1.Range("H6").Value2 = 3

# row 7
$q1.Range("A7").Value2 = 5
$q1.Range("B7").Value2 = '''166023'
$q1.Range("C7").Value2 = '中欧瑞丰灵活配置混合（LOF）A'
$q1.Range("D7").Value2 = '''32.40'
$q1.Range("E7").Value2 = '''85.04'
$q1.Range("F7").Value2 = '''4.44'
$q1.Range("G7").Value2 = '''1.4386'
$q1.Range("H7").Value2 = 4

# row 8
$q1.Range("A8").Value2 = 6
$q1.Range("B8").Value2 = '''007549'
$q1.Range("C8").Value2 = '中泰开阳价值优选灵活配置混合A'
$q1.Range("D8").Value2 = '''20.45'
$q1.Range("E8").Value2 = '''91.42'
$q1.Range("F8").Value2 = '''4.03'
$q1.Range("G8").Value2 = '''0.8241'
$q1.Range("H8").Value2 = 10

# row 9
$q1.Range("A9").Value2 = 7
$q1.Range("B9").Value2 = '''000242'
$q1.Range("C9").Value2 = '景顺长城策略精选'
$q1.Range("D9").Value2 = '''20.09'
$q1.Range("E9").Value2 = '''91.24'
$q1.Range("F9").Value2 = '''3.56'
$q1.Range("G9").Value2 = '''0.7152'
$q1.Range("H9").Value2 = 8

# row 10
$q1.Range("A10").Value2 = 8
$q1.Range("B10").Value2 = '''012001'
$q1.Range("C10").Value2 = '中泰星宇价值成长一年封闭运作混合型证券投资基金A'
$q1.Range("D10").Value2 = '''12.14'
$q1.Range("E10").Value2 = '''91.96'
$q1.Range("F10").Value2 = '''4.37'
$q1.Range("G10").Value2 = '''0.5305'
$q1.Range("H10").Value2 = 9

# row 11
$q1.Range("A11").Value2 = 9
$q1.Range("B11").Value2 = '''013414'
$q1.Range("C11").Value2 = '太平智远三个月定期开放股票'
$q1.Range("D11").Value2 = '''8.69'
$q1.Range("E11").Value2 = '''86.34'
$q1.Range("F11").Value2 = '''5.98'
$q1.Range("G11").Value2 = '''0.5197'
$q1.Range("H11").Value2 = 6

# row 12
$q1.Range("A12").Value2 = 10
$q1.Range("B12").Value2 = '''009791'
$q1.Range("C12").Value2 = '中欧创业板两年定期开放混合C'
$q1.Range("D12").Value2 = '''5.21'
$q1.Range("E12").Value2 = '''99.64'
$q1.Range("F12").Value2 = '''9.93'
$q1.Range("G12").Value2 = '''0.5174'
$q1.Range("H12").Value2 = 1

# row 13
$q1.Range("A13").Value2 = 11
$q1.Range("B13").Value2 = '''011437'
$q1.Range("C13").Value2 = '中泰开阳价值优选灵活配置混合C'
$q1.Range("D13").Value2 = '''5.02'
$q1.Range("E13").Value2 = '''91.42'
$q1.Range("F13").Value2 = '''4.03'
$q1.Range("G13").Value2 = '''0.2023'
$q1.Range("H13").Value2 = 10

# row 14
$q1.Range("A14").Value2 = 12
$q1.Range("B14").Value2 = '''562510'
$q1.Range("C14").Value2 = '华夏中证旅游主题ETF'
$q1.Range("D14").Value2 = '''1.71'
$q1.Range("E14").Value2 = '''98.99'
$q1.Range("F14").Value2 = '''10.70'
$q1.Range("G14").Value2 = '''0.1830'
$q1.Range("H14").Value2 = 3

# row 15
$q1.Range("A15").Value2 = 13
$q1.Range("B15").Value2 = '''005270'
$q1.Range("C15").Value2 = '太平改革红利精选灵活配置混合'
$q1.Range("D15").Value2 = '''1.87'
$q1.Range("E15").Value2 = '''88.32'
$q1.Range("F15").Value2 = '''6.58'
$q1.Range("G15").Value2 = '''0.1230'
$q1.Range("H15").Value2 = 5

# row 16
$q1.Range("A16").Value2 = 14
$q1.Range("B16").Value2 = '''012002'
$q1.Range("C16").Value2 = '中泰星宇价值成长一年封闭运作混合型证券投资基金C'
$q1.Range("D16").Value2 = '''2.13'
$q1.Range("E16").Value2 = '''91.96'
$q1.Range("F16").Value2 = '''4.37'
$q1.Range("G16").Value2 = '''0.0931'
$q1.Range("H16").Value2 = 9

# row 17
$q1.Range("A17").Value2 = 15
$q1.Range("B17").Value2 = '''010896'
$q1.Range("C17").Value2 = '太平价值增长股票A'
$q1.Range("D17").Value2 = '''1.18'
$q1.Range("E17").Value2 = '''83.63'
$q1.Range("F17").Value2 = '''7.07'
$q1.Range("G17").Value2 = '''0.0834'
$q1.Range("H17").Value2 = 5

# row 18
$q1.Range("A18").Value2 = 16
$q1.Range("B18").Value2 = '''010897'
$q1.Range("C18").Value2 = '太平价值增长股票C'
$q1.Range("D18").Value2 = '''1.01'
$q1.Range("E18").Value2 = '''83.63'
$q1.Range("F18").Value2 = '''7.07'
$q1.Range("G18").Value2 = '''0.0714'
$q1.Range("H18").Value2 = 5

# row 19
$q1.Range("A19").Value2 = 17
$q1.Range("B19").Value2 = '''770001'
$q1.Range("C19").Value2 = '德邦优化灵活配置混合'
$q1.Range("D19").Value2 = '''2.49'
$q1.Range("E19").Value2 = '''86.80'
$q1.Range("F19").Value2 = '''2.86'
$q1.Range("G19").Value2 = '''0.0712'
$q1.Range("H19").Value2 = 9

# row 20
$q1.Range("A20").Value2 = 18
$q1.Range("B20").Value2 = '''004740'
$q1.Range("C20").Value2 = '中欧瑞丰灵活配置混合（LOF）C'
$q1.Range("D20").Value2 = '''1.28'
$q1.Range("E20").Value2 = '''85.04'
$q1.Range("F20").Value2 = '''4.44'
$q1.Range("G20").Value2 = '''0.0568'
$q1.Range("H20").Value2 = 4

# row 21
$q1.Range("A21").Value2 = 19
$q1.Range("B21").Value2 = '''290012'
$q1.Range("C21").Value2 = '泰信行业精选灵活配置混合A'
$q1.Range("D21").Value2 = '''0.76'
$q1.Range("E21").Value2 = '''92.62'
$q1.Range("F21").Value2 = '''5.49'
$q1.Range("G21").Value2 = '''0.0417'
$q1.Range("H21").Value2 = 6

# row 22
$q1.Range("A22").Value2 = 20
$q1.Range("B22").Value2 = '''161036'
$q1.Range("C22").Value2 = '富国中证娱乐主题指数增强（LOF）'
$q1.Range("D22").Value2 = '''0.77'
$q1.Range("E22").Value2 = '''93.32'
$q1.Range("F22").Value2 = '''3.44'
$q1.Range("G22").Value2 = '''0.0265'
$q1.Range("H22").Value2 = 7

# row 23
$q1.Range("A23").Value2 = 21
$q1.Range("B23").Value2 = '''003132'
$q1.Range("C23").Value2 = '德邦新回报灵活配置混合'
$q1.Range("D23").Value2 = '''0.55'
$q1.Range("E23").Value2 = '''79.88'
$q1.Range("F23").Value2 = '''4.72'
$q1.Range("G23").Value2 = '''0.0260'
$q1.Range("H23").Value2 = 1

# row 24
$q1.Range("A24").Value2 = 22
$q1.Range("B24").Value2 = '''001448'
$q1.Range("C24").Value2 = '华商双翼平衡混合'
$q1.Range("D24").Value2 = '''0.38'
$q1.Range("E24").Value2 = '''39.74'
$q1.Range("F24").Value2 = '''2.07'
$q1.Range("G24").Value2 = '''0.0079'
$q1.Range("H24").Value2 = 6

# row 25
$q1.Range("A25").Value2 = 23
$q1.Range("B25").Value2 = '''516190'
$q1.Range("C25").Value2 = '华夏中证文娱传媒ETF'
$q1.Range("D25").Value2 = '''0.07'
$q1.Range("E25").Value2 = '''96.81'
$q1.Range("F25").Value2 = '''3.60'
$q1.Range("G25").Value2 = '''0.0025'
$q1.Range("H25").Value2 = 7

# row 26
$q1.Range("A26").Value2 = 24
$q1.Range("B26").Value2 = '''001466'
$q1.Range("C26").Value2 = '华富永鑫灵活配置混合A'
$q1.Range("D26").Value2 = '''0.05'
$q1.Range("E26").Value2 = '''29.62'
$q1.Range("F26").Value2 = '''0.96'
$q1.Range("G26").Value2 = '''0.0005'
$q1.Range("H26").Value2 = 10

# row 27
$q1.Range("A27").Value2 = 25
$q1.Range("B27").Value2 = '''001467'
$q1.Range("C27").Value2 = '华富永鑫灵活配置混合C'
$q1.Range("D27").Value2 = '''0.04'
$q1.Range("E27").Value2 = '''29.62'
$q1.Range("F27").Value2 = '''0.96'
$q1.Range("G27").Value2 = '''0.0004'
$q1.Range("H27").Value2 = 10

# row 28
$q1.Range("A28").Value2 = 26
$q1.Range("B28").Value2 = '''002583'
$q1.Range("C28").Value2 = '泰信行业精选灵活配置混合C'
$q1.Range("D28").Value2 = '''0.00'
$q1.Range("E28").Value2 = '''92.62'
$q1.Range("F28").Value2 = '''5.49'
$q1.Range("G28").Value2 = 0
$q1.Range("H28").Value2 = 6

# ---- Step 4: insert the new 2022-Q1 summary row at the top of the '总计' sheet's data ----
$ws2 = $wb.Worksheets.Item("总计")

# Shift existing data rows (2..6) down by one to make room for the new 2022-Q1 row.
$ws2.Rows.Item(2).Insert()

# The inserted row inherits the row-above's (header's) bold formatting by default;
# clear it so the new data row matches the plain (unstyled) look of the other data rows.
$ws2.Range("A2:D2").ClearFormats()

# Re-apply the bold/bordered "index column" style (used by every other A-column cell
# in this sheet) onto the new A2 cell by copying it from an existing index cell.
$ws2.Range("A3").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("A2").Value2 = 0
$ws2.Range("B2").Value2 = "2022-Q1"
$ws2.Range("C2").Value2 = 27
$ws2.Range("D2").Value2 = 17.87

# Renumber the 0-based index column for the rows that shifted down.
$ws2.Range("A3").Value2 = 1
$ws2.Range("A4").Value2 = 2
$ws2.Range("A5").Value2 = 3
$ws2.Range("A6").Value2 = 4
$ws2.Range("A7").Value2 = 5
